# Fruta / hortaliza, semanal
# The source publishes a new week's worth of rows at the top of the
# per-quality-grade block (rows 131-132) and pushes the existing history
# down by two rows; the two rows that fall off the bottom (old 216/217)
# are re-appended as new rows 218/219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 131; this shifts rows 131-217 down to
# 133-219 (values, formats and styles all move with the rows), which
# automatically reproduces the bulk of the diff and grows the used range
# to A1:T219.
$ws.Range("A131:A132").EntireRow.Insert()

# Populate the first of the two newly-inserted rows (131) with the new
# week's "Especial" grade entry.
$ws.Range("A131").Value = 3
$ws.Range("B131").Value = "Femacal de La Calera"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44567
$ws.Range("D131").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E131").Value = 5
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100101
$ws.Range("H131").Value = "Berries"
$ws.Range("I131").Value = 100112025
$ws.Range("J131").Value = "Frutilla"
$ws.Range("K131").Value = "Sin especificar"
$ws.Range("L131").Value = "Especial"
$ws.Range("M131").Value = 58
$ws.Range("N131").Value = 7000
$ws.Range("O131").Value = 7000
$ws.Range("P131").Value = 7000
$ws.Range("Q131").Value = "$/bandeja 7 kilos"
$ws.Range("R131").Value = "Provincia de Melipilla"
$ws.Range("S131").Value = 1000
$ws.Range("T131").Value = 7

# Populate the second newly-inserted row (132) with the new week's
# "Segunda" grade entry.
$ws.Range("A132").Value = 3
$ws.Range("B132").Value = "Femacal de La Calera"
$ws.Range("C132").Value = "Coquimbo"
$ws.Range("D132").Value = 44567
$ws.Range("D132").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E132").Value = 5
$ws.Range("F132").Value = "Fruta"
$ws.Range("G132").Value = 100101
$ws.Range("H132").Value = "Berries"
$ws.Range("I132").Value = 100112025
$ws.Range("J132").Value = "Frutilla"
$ws.Range("K132").Value = "Sin especificar"
$ws.Range("L132").Value = "Segunda"
$ws.Range("M132").Value = 48
$ws.Range("N132").Value = 4000
$ws.Range("O132").Value = 4000
$ws.Range("P132").Value = 4000
$ws.Range("Q132").Value = "$/bandeja 7 kilos"
$ws.Range("R132").Value = "Provincia de Melipilla"
$ws.Range("S132").Value = 571
$ws.Range("T132").Value = 7
